$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Datas da campaña de Constelación de Perseo 2022:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Datas da campaña de 2022 que usan Constelación de Perseo:",
    2
)
